$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the changed cells keep their text (string) representation
# instead of being auto-converted to numbers/percentages by Excel.
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "287.57"
$ws.Range("E2").Value = "-10.10%"
$ws.Range("G2").Value = "18"
$ws.Range("D3").Value = "40.14"
$ws.Range("E3").Value = "-3.20%"
$ws.Range("G3").Value = "18"
$ws.Range("D4").Value = "5.038"
$ws.Range("E4").Value = "-3.86%"
$ws.Range("G4").Value = "18"
$ws.Range("D5").Value = "0.07275"
$ws.Range("E5").Value = "-5.87%"
$ws.Range("G5").Value = "18"
$ws.Range("D6").Value = "4.283"
$ws.Range("E6").Value = "-1.32%"
$ws.Range("G6").Value = "18"
$ws.Range("D7").Value = "1.511"
$ws.Range("E7").Value = "-10.37%"
$ws.Range("G7").Value = "18"
$ws.Range("D8").Value = "0.9137"
$ws.Range("E8").Value = "-3.40%"
$ws.Range("G8").Value = "18"
$ws.Range("D9").Value = "0.1196"
$ws.Range("E9").Value = "-3.56%"
$ws.Range("G9").Value = "18"
$ws.Range("E10").Value = "-7.25%"
$ws.Range("G10").Value = "18"
$ws.Range("D11").Value = "0.08519"
$ws.Range("E11").Value = "-7.37%"
$ws.Range("G11").Value = "18"
$ws.Range("D12").Value = "0.04155"
$ws.Range("E12").Value = "-3.96%"
$ws.Range("G12").Value = "18"
$ws.Range("D13").Value = "0.1054"
$ws.Range("E13").Value = "0.41%"
$ws.Range("G13").Value = "18"
$ws.Range("D14").Value = "0.001278"
$ws.Range("E14").Value = "-0.69%"
$ws.Range("G14").Value = "18"
$ws.Range("D15").Value = "0.005993"
$ws.Range("E15").Value = "0.30%"
$ws.Range("G15").Value = "18"
$ws.Range("D16").Value = "3.396"
$ws.Range("E16").Value = "1.67%"
$ws.Range("G16").Value = "18"
$ws.Range("G17").Value = "18"
$ws.Range("D18").Value = "0.3262"
$ws.Range("E18").Value = "-2.91%"
$ws.Range("G18").Value = "18"
$ws.Range("D19").Value = "7.786"
$ws.Range("E19").Value = "1.52%"
$ws.Range("G19").Value = "18"
$ws.Range("D20").Value = "0.1353"
$ws.Range("E20").Value = "-0.07%"
$ws.Range("G20").Value = "18"
$ws.Range("D21").Value = "0.2890"
$ws.Range("E21").Value = "2.30%"
$ws.Range("G21").Value = "18"
$ws.Range("D22").Value = "0.03846"
$ws.Range("E22").Value = "-4.81%"
$ws.Range("G22").Value = "18"
$ws.Range("D23").Value = "0.001269"
$ws.Range("E23").Value = "0.29%"
$ws.Range("G23").Value = "18"
$ws.Range("D24").Value = "0.003804"
$ws.Range("E24").Value = "-7.77%"
$ws.Range("G24").Value = "18"
$ws.Range("D25").Value = "0.0001282"
$ws.Range("E25").Value = "0.89%"
$ws.Range("G25").Value = "18"
$ws.Range("D26").Value = "0.0003729"
$ws.Range("G26").Value = "18"
$ws.Range("G27").Value = "18"
$ws.Range("G28").Value = "18"
$ws.Range("G29").Value = "18"
$ws.Range("G30").Value = "18"
$ws.Range("G31").Value = "18"
$ws.Range("G32").Value = "18"
$ws.Range("G33").Value = "18"
$ws.Range("G34").Value = "18"
$ws.Range("G35").Value = "18"
$ws.Range("G36").Value = "18"
$ws.Range("G37").Value = "18"
$ws.Range("D38").Value = "0.02294"
$ws.Range("E38").Value = "-9.84%"
$ws.Range("G38").Value = "18"
$ws.Range("D39").Value = "0.04933"
$ws.Range("E39").Value = "-7.71%"
$ws.Range("G39").Value = "18"
$ws.Range("D40").Value = "0.006893"
$ws.Range("E40").Value = "246.03%"
$ws.Range("G40").Value = "18"
$ws.Range("D41").Value = "0.007713"
$ws.Range("E41").Value = "-0.46%"
$ws.Range("G41").Value = "18"
$ws.Range("D42").Value = "0.1267"
$ws.Range("E42").Value = "-3.91%"
$ws.Range("G42").Value = "18"
$ws.Range("D43").Value = "0.007378"
$ws.Range("E43").Value = "0.29%"
$ws.Range("G43").Value = "18"
$ws.Range("D44").Value = "0.007690"
$ws.Range("E44").Value = "-8.00%"
$ws.Range("G44").Value = "18"
$ws.Range("D45").Value = "0.3110"
$ws.Range("E45").Value = "-10.13%"
$ws.Range("G45").Value = "18"
$ws.Range("D46").Value = "0.00006389"
$ws.Range("E46").Value = "-4.62%"
$ws.Range("G46").Value = "18"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").Value = "0.09%"
$ws.Range("G47").Value = "18"
$ws.Range("E48").Value = "23.73%"
$ws.Range("G48").Value = "18"
$ws.Range("E49").Value = "-0.06%"
$ws.Range("G49").Value = "18"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").Value = "0.09%"
$ws.Range("G50").Value = "18"
$ws.Range("D51").Value = "0.0002004"
$ws.Range("E51").Value = "0.09%"
$ws.Range("G51").Value = "18"
